# Auto-generated update of market-price derived columns (H-N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Goblin Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2427.375
$ws.Range("J17").Value = 2427.375
$ws.Range("L17").Value = 7282.125
$ws.Range("N17").Value = -7618.125
$ws.Range("H53").Value = 177.52942
$ws.Range("I53").Value = 92.333336
$ws.Range("K53").Value = 92.333336
$ws.Range("M53").Value = 544.666664
$ws.Range("H80").Value = 3015.7
$ws.Range("I80").Value = 718.6667
$ws.Range("K80").Value = 2156.0001
$ws.Range("M80").Value = -1158.0001
$ws.Range("H83").Value = 3015.7
$ws.Range("I83").Value = 718.6667
$ws.Range("K83").Value = 6468.0003
$ws.Range("M83").Value = -1476.0003
$ws.Range("H137").Value = 1813.9131
$ws.Range("I137").Value = 1886
$ws.Range("K137").Value = 5658
$ws.Range("M137").Value = -3108
$ws.Range("H141").Value = 9110.883
$ws.Range("I141").Value = 8824
$ws.Range("K141").Value = 26472
$ws.Range("M141").Value = -21292
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9704.294
$ws.Range("I32").Value = 9704.294
$ws.Range("K32").Value = 9704.294
$ws.Range("M32").Value = -9417.294
$ws.Range("H45").Value = 2374.0557
$ws.Range("I45").Value = 1896
$ws.Range("K45").Value = 1896
$ws.Range("M45").Value = -1519
$ws.Range("H61").Value = 5911.4165
$ws.Range("I61").Value = 5876.409
$ws.Range("K61").Value = 5876.409
$ws.Range("M61").Value = -5664.409
$ws.Range("H74").Value = 2787.9375
$ws.Range("I74").Value = 2720
$ws.Range("K74").Value = 2720
$ws.Range("M74").Value = -1846
$ws.Range("H77").Value = 2787.9375
$ws.Range("I77").Value = 2720
$ws.Range("K77").Value = 13600
$ws.Range("M77").Value = -9232
$ws.Range("H102").Value = 4346.647
$ws.Range("I102").Value = 1263.091
$ws.Range("J102").Value = 9999.833000000001
$ws.Range("K102").Value = 1263.091
$ws.Range("L102").Value = 9999.833000000001
$ws.Range("M102").Value = 358.9090000000001
$ws.Range("N102").Value = -13243.833
$ws.Range("H110").Value = 2033.64
$ws.Range("I110").Value = 2033.64
$ws.Range("K110").Value = 2033.64
$ws.Range("M110").Value = 11.3599999999999
$ws.Range("H122").Value = 5557494.5
$ws.Range("I122").Value = 7409051.5
$ws.Range("K122").Value = 22227154.5
$ws.Range("M122").Value = -22224704.5
$ws.Range("H132").Value = 7158.8945
$ws.Range("I132").Value = 7445.5
$ws.Range("K132").Value = 22336.5
$ws.Range("M132").Value = -19806.5
$ws.Range("H135").Value = 49827.5
$ws.Range("J135").Value = 49827.5
$ws.Range("L135").Value = 49827.5
$ws.Range("N135").Value = -59967.5
$ws.Range("H136").Value = 5911.4165
$ws.Range("I136").Value = 5876.409
$ws.Range("K136").Value = 17629.227
$ws.Range("M136").Value = -15079.227
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1576.0385
$ws.Range("I20").Value = 1475.4117
$ws.Range("J20").Value = 1766.1111
$ws.Range("K20").Value = 1475.4117
$ws.Range("L20").Value = 1766.1111
$ws.Range("M20").Value = -1228.4117
$ws.Range("N20").Value = -2260.1111
$ws.Range("H30").Value = 1610.8
$ws.Range("I30").Value = 10
$ws.Range("K30").Value = 10
$ws.Range("M30").Value = 115
$ws.Range("H40").Value = 38851.75
$ws.Range("I40").Value = 35407
$ws.Range("K40").Value = 35407
$ws.Range("M40").Value = -35142
$ws.Range("H75").Value = 10302.2
$ws.Range("I75").Value = 10302.2
$ws.Range("K75").Value = 10302.2
$ws.Range("M75").Value = -9366.200000000001
$ws.Range("H78").Value = 10302.2
$ws.Range("I78").Value = 10302.2
$ws.Range("K78").Value = 30906.6
$ws.Range("M78").Value = -26226.6
$ws.Range("H86").Value = 18521416
$ws.Range("I86").Value = 1702.5
$ws.Range("J86").Value = 71434890
$ws.Range("K86").Value = 1702.5
$ws.Range("L86").Value = 71434890
$ws.Range("M86").Value = -579.5
$ws.Range("N86").Value = -71437136
$ws.Range("H89").Value = 18521416
$ws.Range("I89").Value = 1702.5
$ws.Range("J89").Value = 71434890
$ws.Range("K89").Value = 8512.5
$ws.Range("L89").Value = 357174450
$ws.Range("M89").Value = -2896.5
$ws.Range("N89").Value = -357185682
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2007.579
$ws.Range("I134").Value = 2009.4375
$ws.Range("K134").Value = 6028.3125
$ws.Range("M134").Value = -3493.3125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 138.46666
$ws.Range("I12").Value = 115.5
$ws.Range("J12").Value = 153.77777
$ws.Range("K12").Value = 346.5
$ws.Range("L12").Value = 461.33331
$ws.Range("M12").Value = -173.5
$ws.Range("N12").Value = -807.33331
$ws.Range("H23").Value = 3181.7856
$ws.Range("I23").Value = 6757
$ws.Range("J23").Value = 2206.7273
$ws.Range("K23").Value = 20271
$ws.Range("L23").Value = 6620.1819
$ws.Range("M23").Value = -20036
$ws.Range("N23").Value = -7090.1819
$ws.Range("H33").Value = 339.46155
$ws.Range("J33").Value = 349.33334
$ws.Range("L33").Value = 2096.00004
$ws.Range("N33").Value = -2662.00004
$ws.Range("H64").Value = 1136.3334
$ws.Range("I64").Value = 1136.3334
$ws.Range("K64").Value = 3409.0002
$ws.Range("M64").Value = -3139.0002
$ws.Range("H67").Value = 1136.3334
$ws.Range("I67").Value = 1136.3334
$ws.Range("K67").Value = 3409.0002
$ws.Range("M67").Value = -2473.0002
$ws.Range("H121").Value = 564.75
$ws.Range("I121").Value = 419.66666
$ws.Range("K121").Value = 1258.99998
$ws.Range("M121").Value = 51.00001999999995
$ws.Range("H131").Value = 4765111.5
$ws.Range("I131").Value = 1333.3334
$ws.Range("J131").Value = 6064323.5
$ws.Range("K131").Value = 4000.0002
$ws.Range("L131").Value = 18192970.5
$ws.Range("M131").Value = 1039.9998
$ws.Range("N131").Value = -18203050.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15882671
$ws.Range("I70").Value = 25649730
$ws.Range("J70").Value = 11198.875
$ws.Range("K70").Value = 25649730
$ws.Range("L70").Value = 11198.875
$ws.Range("M70").Value = -25649460
$ws.Range("N70").Value = -11738.875
$ws.Range("H73").Value = 15882671
$ws.Range("I73").Value = 25649730
$ws.Range("J73").Value = 11198.875
$ws.Range("K73").Value = 25649730
$ws.Range("L73").Value = 11198.875
$ws.Range("M73").Value = -25648794
$ws.Range("N73").Value = -13070.875
$ws.Range("H80").Value = 4861.7607
$ws.Range("I80").Value = 3076.8262
$ws.Range("J80").Value = 6646.696
$ws.Range("K80").Value = 3076.8262
$ws.Range("L80").Value = 6646.696
$ws.Range("M80").Value = -2078.8262
$ws.Range("N80").Value = -8642.696
$ws.Range("H83").Value = 4861.7607
$ws.Range("I83").Value = 3076.8262
$ws.Range("J83").Value = 6646.696
$ws.Range("K83").Value = 15384.131
$ws.Range("L83").Value = 33233.48
$ws.Range("M83").Value = -10392.131
$ws.Range("N83").Value = -43217.48
$ws.Range("H122").Value = 7040.864
$ws.Range("I122").Value = 7121.1055
$ws.Range("K122").Value = 21363.3165
$ws.Range("M122").Value = -18913.3165
$ws.Range("H126").Value = 4316.3335
$ws.Range("I126").Value = 2974.5
$ws.Range("K126").Value = 8923.5
$ws.Range("M126").Value = -6453.5
$ws.Range("H132").Value = 2587.1052
$ws.Range("I132").Value = 2213.8708
$ws.Range("K132").Value = 6641.6124
$ws.Range("M132").Value = -4111.6124
$ws.Range("H134").Value = 78886.8
$ws.Range("J134").Value = 78886.8
$ws.Range("L134").Value = 236660.4
$ws.Range("N134").Value = -241730.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2751.5715
$ws.Range("I16").Value = 768.86664
$ws.Range("J16").Value = 7708.3335
$ws.Range("K16").Value = 768.86664
$ws.Range("L16").Value = 7708.3335
$ws.Range("M16").Value = -598.86664
$ws.Range("N16").Value = -8048.3335
$ws.Range("H46").Value = 2806.1875
$ws.Range("I46").Value = 499
$ws.Range("K46").Value = 499
$ws.Range("M46").Value = -311
$ws.Range("H82").Value = 2340.182
$ws.Range("I82").Value = 392
$ws.Range("J82").Value = 3070.75
$ws.Range("K82").Value = 392
$ws.Range("L82").Value = 3070.75
$ws.Range("M82").Value = -31
$ws.Range("N82").Value = -3792.75
$ws.Range("H85").Value = 2340.182
$ws.Range("I85").Value = 392
$ws.Range("J85").Value = 3070.75
$ws.Range("K85").Value = 392
$ws.Range("L85").Value = 3070.75
$ws.Range("M85").Value = 856
$ws.Range("N85").Value = -5566.75
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49972
$ws.Range("J46").Value = 49972
$ws.Range("L46").Value = 49972
$ws.Range("N46").Value = -50434
$ws.Range("H81").Value = 2154.6667
$ws.Range("I81").Value = 1607.75
$ws.Range("K81").Value = 3215.5
$ws.Range("M81").Value = -2154.5
$ws.Range("H84").Value = 2154.6667
$ws.Range("I84").Value = 1607.75
$ws.Range("K84").Value = 16077.5
$ws.Range("M84").Value = -10773.5
$ws.Range("H92").Value = 19999.5
$ws.Range("J92").Value = 19999.5
$ws.Range("L92").Value = 19999.5
$ws.Range("N92").Value = -24991.5
$ws.Range("H100").Value = 1202.4445
$ws.Range("I100").Value = 834
$ws.Range("J100").Value = 1386.6666
$ws.Range("K100").Value = 1668
$ws.Range("L100").Value = 2773.3332
$ws.Range("M100").Value = -1127
$ws.Range("N100").Value = -3855.3332
$ws.Range("H134").Value = 49972
$ws.Range("J134").Value = 49972
$ws.Range("L134").Value = 149916
$ws.Range("N134").Value = -154986
